$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column L: 2021 header (row 4) and its value 269 (row 5)
$ws.Range("L4").Value = 2021
$ws.Range("L5").Value = 269

# Extend the thin bottom-border formatting row (row 3) into the new column,
# matching the existing K column style (also picks up D:K's header/value
# cell styles for L4/L5 without clobbering the values just written).
$ws.Range("K3:K5").Copy()
$ws.Range("L3:L5").PasteSpecial(-4122)

# Match the author's final selection recorded in the workbook
$ws.Range("N3").Select()
